# Finish the "submitting po based invoice manually" section:
#  - Refresh the sample invoice rows on the POBasedInvoice sheet (new
#    invoice numbers / base amounts / IGST, and make Quantity mirror the
#    Base Amount instead of a flat literal 1).
#  - Correct the BADashboardPage sample "to state" value and widen its
#    second column.
#  - Leave BADashboardPage as the active/selected sheet, matching the
#    state the workbook was saved in.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# POBasedInvoice sheet: rows 2-9 (Invoice Number / Base Amount / IGST /
# Quantity). Base Amount, IGST and Quantity are stored as text, not
# numbers, in this workbook, so values that look numeric are entered
# with a leading apostrophe to force text just like the original data.
# ---------------------------------------------------------------------
$po = $wb.Worksheets.Item("POBasedInvoice")

$invoiceRows = @(
    @{ Row = 2;  Invoice = "TESTINV24257"; Base = "7"; Igst = "1.26" },
    @{ Row = 3;  Invoice = "TESTINV70287"; Base = "8"; Igst = "1.44" },
    @{ Row = 4;  Invoice = "TESTINV07707"; Base = "9"; Igst = "1.62" },
    @{ Row = 5;  Invoice = "TESTINV49593"; Base = "9"; Igst = "1.62" },
    @{ Row = 6;  Invoice = "TESTINV87128"; Base = "3"; Igst = "0.54" },
    @{ Row = 7;  Invoice = "TESTINV79234"; Base = "8"; Igst = "1.44" },
    @{ Row = 8;  Invoice = "TESTINV09626"; Base = "3"; Igst = "0.54" },
    @{ Row = 9;  Invoice = "TESTINV34765"; Base = "7"; Igst = "1.26" }
)

foreach ($r in $invoiceRows) {
    $po.Range("A" + $r.Row).Value = $r.Invoice
    $po.Range("B" + $r.Row).Value = "'" + $r.Base
    $po.Range("C" + $r.Row).Value = "'" + $r.Igst
    $po.Range("J" + $r.Row).Value = "'" + $r.Base
}

# ---------------------------------------------------------------------
# BADashboardPage sheet: fix the sample "to state" value and widen
# column B to fit it.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("BADashboardPage")
$dash.Range("B2").Value = "TRIPURA"
$dash.Columns.Item(2).ColumnWidth = 14.6640625

# Leave BADashboardPage as the selected/active sheet.
$dash.Activate()
